# Weekly data refresh: insert one new price record for
# "Vega Modelo de Temuco - Pomelo" as row 95, pushing the existing
# rows 95-141 down to 96-142 (dimension grows from T141 to T142).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 95; everything currently at/after row 95 shifts down one.
$ws.Rows(95).Insert()

# Populate the new row 95 with the new weekly record.
$ws.Range("A95").Value = 10
$ws.Range("B95").Value = "Vega Modelo de Temuco"
$ws.Range("C95").Value = "La Araucanía"
$ws.Range("D95").Value = 44466
$ws.Range("E95").Value = 9
$ws.Range("F95").Value = "Fruta"
$ws.Range("G95").Value = 100102
$ws.Range("H95").Value = "Cítricos"
$ws.Range("I95").Value = 100102006
$ws.Range("J95").Value = "Pomelo"
$ws.Range("K95").Value = "Start Ruby"
$ws.Range("L95").Value = "Especial"
$ws.Range("M95").Value = 100
$ws.Range("N95").Value = 10000
$ws.Range("O95").Value = 10000
$ws.Range("P95").Value = 10000
$ws.Range("Q95").Value = "$/bandeja 15 kilos granel"
$ws.Range("R95").Value = "Región de O'Higgins"
$ws.Range("S95").Value = 667
$ws.Range("T95").Value = 15
